# Updates cryptos list values (Coin/Link/Price/Volume(1h)) to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($Cell, $Text)
    $range = $ws.Range($Cell)
    $range.Value = $Text
    # Re-apply the default "Normal" style so plain numeric-looking
    # strings (forced to text via a leading apostrophe) keep the
    # same formatting/style as the rest of the worksheet.
    $range.Style = "Normal"
}

Set-CellText "D2" "27.400.40"
Set-CellText "D3" "1.641.39"
Set-CellText "E3" "  -1.52%  "
Set-CellText "E4" "  +0.01%  "
Set-CellText "D5" "'211.92"
Set-CellText "E5" "  -1.51%  "
Set-CellText "E6" "  +4.24%  "
Set-CellText "E7" "  +0.00%  "
Set-CellText "D8" "'23.18"
Set-CellText "E8" "  -1.57%  "
Set-CellText "E9" "  -2.22%  "
Set-CellText "E10" "  -2.01%  "
Set-CellText "D11" "'0.0891"
Set-CellText "E11" "  +1.29%  "
Set-CellText "D12" "1.873.91"
Set-CellText "E12" "  -1.50%  "
Set-CellText "D13" "1.651.01"
Set-CellText "E13" "  -0.71%  "
Set-CellText "D14" "'4.02"
Set-CellText "E14" "  -3.21%  "
Set-CellText "D15" "'0.558"
Set-CellText "E15" "  +0.71%  "
Set-CellText "D16" "'64.22"
Set-CellText "E16" "  -3.12%  "
Set-CellText "D17" "27.376.77"
Set-CellText "E17" "  -0.78%  "
Set-CellText "D18" "'227.70"
Set-CellText "E18" "  -9.25%  "
Set-CellText "E19" "  -1.83%  "
Set-CellText "D20" "'7.46"
Set-CellText "E20" "  -1.09%  "
Set-CellText "D21" "'1.00"
Set-CellText "E21" "  +0.02%  "
Set-CellText "E22" "  -4.31%  "
Set-CellText "D23" "'9.28"
Set-CellText "E23" "  -0.24%  "
Set-CellText "E24" "  +0.25%  "
Set-CellText "D25" "'147.61"
Set-CellText "E25" "  +0.75%  "
Set-CellText "D26" "'0.115"
Set-CellText "E26" "  +2.44%  "
Set-CellText "E27" "  -3.03%  "
Set-CellText "E28" "  +0.05%  "
Set-CellText "D29" "'15.51"
Set-CellText "E29" "  -6.39%  "
Set-CellText "E30" "  -4.91%  "
Set-CellText "E31" "  -4.18%  "
Set-CellText "E32" "  -2.72%  "
Set-CellText "E33" "  -0.62%  "
Set-CellText "D34" "1.397.48"
Set-CellText "E34" "  -5.21%  "
Set-CellText "E35" "  -1.16%  "
Set-CellText "E36" "  -0.28%  "
Set-CellText "E37" "  -3.06%  "
Set-CellText "E38" "  -7.02%  "
Set-CellText "D39" "'0.0166"
Set-CellText "E39" "  -3.12%  "
Set-CellText "E40" "  -0.31%  "
Set-CellText "D41" "'1.00"
Set-CellText "E41" "  +0.01%  "
Set-CellText "B42" "FraxShare"
Set-CellText "C42" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-CellText "D42" "'5.47"
Set-CellText "E42" "  +0.73%  "
Set-CellText "B43" "MXToken"
Set-CellText "C43" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-CellText "D43" "'2.21"
Set-CellText "E43" "  +0.05%  "
Set-CellText "B44" "TrustWalletToken"
Set-CellText "C44" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-CellText "D44" "'0.788"
Set-CellText "E44" "  -0.38%  "
Set-CellText "B45" "Aave"
Set-CellText "C45" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-CellText "D45" "'64.20"
Set-CellText "E45" "  -7.89%  "
Set-CellText "B46" "RocketPoolETH"
Set-CellText "C46" "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-CellText "D46" "1.784.41"
Set-CellText "E46" "  -1.41%  "
Set-CellText "B47" "RenderToken"
Set-CellText "C47" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-CellText "D47" "'1.64"
Set-CellText "E47" "  -3.69%  "
Set-CellText "B48" "Quant"
Set-CellText "C48" "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-CellText "D48" "'87.23"
Set-CellText "E48" "  -2.52%  "
Set-CellText "B49" "BabyDogeCoin"
Set-CellText "C49" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-CellText "D49" "0.0₆0105"
Set-CellText "E49" "  -3.66%  "
Set-CellText "B50" "Algorand"
Set-CellText "C50" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-CellText "D50" "'0.0982"
Set-CellText "E50" "  -3.56%  "
Set-CellText "B51" "EnergySwap"
Set-CellText "C51" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-CellText "D51" "'7.61"
Set-CellText "E51" "  -3.90%  "
